# "did the Min stack" - add a new row (17) to the LeetCode summary sheet
# for "155. Min Stack" under the "Stack" category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
$ws.Range("A17").Value = "Stack "
$ws.Range("B17").Value = "155. Min Stack"
$ws.Range("C17").Value = 'We have to design a stack with operations push(), pop(), top() & getMin(), all having time O(1),,, the first 3 are straightforward, but for getMin() brute force is O(n),, for getMin we will make a stack<Node> where node has int num (curr num) & int min (this is the minimum num from curr num to end/bottom of the stack) , while pushing you set Node.num = inputNum & Node.min = min(inputNum, getMin()),, ie u set min = min(input, previous minimum in stack)'

# --- Formatting: match the look of the other rows (thin border all
# round, left/center aligned, wrapped text; column B highlighted
# yellow like the rest of the "Question" column) ----------------------
$newRow = $ws.Range("A17:C17")
$newRow.HorizontalAlignment = -4131   # xlLeft
$newRow.VerticalAlignment = -4108     # xlCenter
$newRow.WrapText = $true
$newRow.Borders.LineStyle = 1
$newRow.Borders.Weight = 2

$ws.Range("B17").Interior.Color = 65535   # yellow fill, same as other Question cells

# Row height matches the other multi-line rows in the sheet
$ws.Rows(17).RowHeight = 57.6

# --- Selection / view ---------------------------------------------------
$ws.Range("E17:E19").Select() | Out-Null
